# Realestate Update resale numbers 2025-02-06 22:16
# Append a new data row (row 55) to the CityResaleNum sheet, mirroring the
# layout of the existing rows (Date, Time, Weekday, Week as text; the city
# columns as numbers).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 55

# Date column: force text first so "2025-02-06" is stored literally instead
# of being reinterpreted as a date serial, then drop back to the default
# (unstyled) cell style so no extra number format sticks to the cell.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "2025-02-06"
$ws.Cells.Item($row, 1).Style = "Normal"

# Time / Weekday columns: plain text, not number/date-like, stores as text
# without any extra formatting.
$ws.Cells.Item($row, 2).Value = "22:16:08"
$ws.Cells.Item($row, 3).Value = "Thursday"

# Week column: force text so "05" keeps its leading zero instead of being
# coerced to the number 5, then reset the style the same way as column A.
$ws.Cells.Item($row, 4).NumberFormat = "@"
$ws.Cells.Item($row, 4).Value = "05"
$ws.Cells.Item($row, 4).Style = "Normal"

# Numeric city columns (E..T)
$ws.Cells.Item($row, 5).Value = 125849
$ws.Cells.Item($row, 6).Value = 141733
$ws.Cells.Item($row, 7).Value = 167750
$ws.Cells.Item($row, 8).Value = 158056
$ws.Cells.Item($row, 9).Value = -1
$ws.Cells.Item($row, 10).Value = 142815
$ws.Cells.Item($row, 11).Value = -1
$ws.Cells.Item($row, 12).Value = -1
$ws.Cells.Item($row, 13).Value = 191286
$ws.Cells.Item($row, 14).Value = 115137
$ws.Cells.Item($row, 15).Value = 44769
$ws.Cells.Item($row, 16).Value = 28281
$ws.Cells.Item($row, 17).Value = 63628
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 41012
$ws.Cells.Item($row, 20).Value = -1
